$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "FLOPADEX 8 MG 30 CAPSULES" product row (row 27). Excel shifts every
# row below it up by one (incl. merged cells, shared-string refs, totals/footer rows).
$ws.Rows.Item(27).Delete()

# The "م" (serial number) column is a manually maintained 1..N sequence, independent
# of the product shift above, so restore it to the correct fixed values per row.
$serials = @(24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63)
for ($i = 0; $i -lt $serials.Length; $i++) {
    $ws.Cells.Item(27 + $i, 1).Value = $serials[$i]
}

# Row heights likewise stay pinned to their original row position (alternating pattern),
# not carried along with the shifted product data.
$heights = @(25.5,25.5,24.75,25.5,24.75,25.5,25.5,24.75,25.5,24.75,25.5,25.5,24.75,25.5,24.75,25.5,25.5,24.75,25.5,24.75,25.5,25.5,24.75,25.5,24.75,25.5,25.5,24.75,25.5,24.75,25.5,25.5,24.75,25.5,24.75,25.5,25.5,24.75,25.5,24.75)
for ($i = 0; $i -lt $heights.Length; $i++) {
    $ws.Rows.Item(27 + $i).RowHeight = $heights[$i]
}

# The grand-total cell is a literal (not a live formula), so correct it for the
# removed row's contribution.
$ws.Range("K67").Value = 3176.1700000000001

# Final total/footer rows pick up their natural (re-fitted) heights.
$ws.Rows.Item(67).RowHeight = 26.25
$ws.Rows.Item(68).RowHeight = 16.5
